$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append new row 12 ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Cells.Item(12, 1).Value = "Demo inplannen"
$ws.Cells.Item(12, 2).Value = "klantenservice@testbedrijf123.nl"
$ws.Cells.Item(12, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Cells.Item(12, 4).Value = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item(12, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Cells.Item(12, 6).Value = "2025-08-13 21:44:48"
$ws.Cells.Item(12, 7).Value = "Nee"
$ws.Cells.Item(12, 8).Value = "Ja"
$ws.Cells.Item(12, 9).Value = "Nee"
$ws.Cells.Item(12, 10).Value = "Nee"

# --- Extend conditional formatting ranges to include the new row (row 12) ---
$ws.Range("D2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D12"))
$ws.Range("G2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G12"))
$ws.Range("H2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H12"))
$ws.Range("I2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I12"))
$ws.Range("J2").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J12"))

# --- "Dashboard" sheet: bump the count for this category ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 11
